# Cotações atualizadas - 2025-10-15
# Adds a new row (41) with the day's quotes, mirroring the formatting of the
# previous row (40): column A keeps the date/number style, columns B-E hold
# the quote values as text (comma decimal separator), matching the rest of
# the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New date (serial 45945 = 2025-10-15) and quotes for the four funds.
$ws.Range("A41").Value = 45945
$ws.Range("B41").Value = "21,6234"
$ws.Range("C41").Value = "15,4836"
$ws.Range("D41").Value = "15,3508"
$ws.Range("E41").Value = "15,3508"

# Match the date-cell's number format to the one used by the existing rows
# (e.g. A40), so the new row stays visually/structurally consistent.
$ws.Range("A41").NumberFormat = $ws.Range("A40").NumberFormat
